# Update "want-to-go" counts (column F) and a couple of lowest-price
# values (column G) to the refreshed numbers from the latest scrape,
# mirroring the change across the per-category sheets (展览, 演出,
# 本地生活) and the combined "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item(1)   # 展览
$sheetShow       = $wb.Worksheets.Item(2)   # 演出
$sheetLocalLife  = $wb.Worksheets.Item(3)   # 本地生活
$sheetAll        = $wb.Worksheets.Item(4)   # 全部类型

function Set-F {
    param($ws, [int]$row, [double]$value)
    $ws.Cells.Item($row, 6).Value = $value   # column F
}

function Set-G {
    param($ws, [int]$row, [double]$value)
    $ws.Cells.Item($row, 7).Value = $value   # column G
}

# ---- 展览 (sheet 1) ----
Set-F $sheetExhibition 3  1204
Set-F $sheetExhibition 4  1155
Set-F $sheetExhibition 6  26
Set-F $sheetExhibition 7  1740
Set-F $sheetExhibition 8  429
Set-F $sheetExhibition 11 322
Set-F $sheetExhibition 12 278
Set-F $sheetExhibition 13 1662
Set-F $sheetExhibition 14 319
Set-F $sheetExhibition 15 1371
Set-F $sheetExhibition 16 771
Set-F $sheetExhibition 17 311
Set-F $sheetExhibition 19 12556
Set-F $sheetExhibition 20 12597
Set-F $sheetExhibition 21 931
Set-F $sheetExhibition 24 291
Set-G $sheetExhibition 24 60
Set-F $sheetExhibition 26 460
Set-F $sheetExhibition 27 1962
Set-F $sheetExhibition 28 19
Set-F $sheetExhibition 29 14

# ---- 演出 (sheet 2) ----
Set-F $sheetShow 5  67
Set-F $sheetShow 7  4
Set-F $sheetShow 9  46
Set-F $sheetShow 10 63

# ---- 本地生活 (sheet 3) ----
Set-F $sheetLocalLife 3 144

# ---- 全部类型 (sheet 4) ----
Set-F $sheetAll 4  1204
Set-F $sheetAll 5  1155
Set-F $sheetAll 7  144
Set-F $sheetAll 8  26
Set-F $sheetAll 9  1740
Set-F $sheetAll 10 429
Set-F $sheetAll 15 322
Set-F $sheetAll 17 278
Set-F $sheetAll 18 1662
Set-F $sheetAll 19 319
Set-F $sheetAll 20 1371
Set-F $sheetAll 21 771
Set-F $sheetAll 22 311
Set-F $sheetAll 23 67
Set-F $sheetAll 25 12556
Set-F $sheetAll 26 12597
Set-F $sheetAll 27 931
Set-F $sheetAll 30 291
Set-G $sheetAll 30 60
Set-F $sheetAll 32 460
Set-F $sheetAll 34 4
Set-F $sheetAll 35 1962
Set-F $sheetAll 36 19
Set-F $sheetAll 38 14
Set-F $sheetAll 39 46
Set-F $sheetAll 42 63
